$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns affected by the weekly rotation: D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# S (Precio $/Kg)
$cols = @("D", "M", "N", "O", "P", "S")

# Capture the original values for rows 3, 4 and 5 before overwriting
# anything, since the update rotates data between these rows:
#   new row3 <- old row4
#   new row4 <- old row5
#   new row5 <- old row3
$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("$col`3").Value2
    $row4[$col] = $ws.Range("$col`4").Value2
    $row5[$col] = $ws.Range("$col`5").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`3").Value2 = $row4[$col]
    $ws.Range("$col`4").Value2 = $row5[$col]
    $ws.Range("$col`5").Value2 = $row3[$col]
}
